$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Edit Repayment Schedule")
$ws.Select()

# The sheet gains a new "waittopageload1" step (same shape as the existing
# "waittopageload" step in row 3) right before the "clickonsubmit" row, so
# insert a blank row at 6 and shift rows 6-12 down to 7-13.
$ws.Rows.Item(6).Insert()

# Row 3's "B" cell (the 2000 value) carries a style (date-ish number format)
# that differs from the plain "fill" style Excel gives a freshly inserted
# row, so copy just that formatting onto the new B6 before writing values.
$ws.Range("B3").Copy()
$ws.Range("B6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A6").Value = "waittopageload1"
$ws.Range("B6").Value = 2000

$ws.Range("A6:B6").Select()
